$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Fill in the previously date-only placeholder row 31 with real data
$ws.Range("B31").Value = "阅读"
$ws.Range("C31").Value = "为什么睡觉"
$ws.Range("D31").Value = 30

# 2. Update the todo note on row 22 (drop the finished "热菜不粘锅" item)
$ws.Range("E22").Value = "todo: swish"

# 3. Remove the now-unneeded extra date-only placeholder row (old row 32),
#    then fully clear the row that shifts up into its place so it doesn't
#    linger as an empty formatted row.
$ws.Rows("32:32").Delete()
$ws.Range("A32").Clear()

# 4. Add the new day's entries starting at row 33
$ws.Range("A33").Value = 44949
$ws.Range("B33").Value = "起床"
$ws.Range("D33").Value = 10

$ws.Range("A34").Value = 44949
$ws.Range("B34").Value = "吃饭"
$ws.Range("D34").Value = 15

$ws.Range("A35").Value = 44949
$ws.Range("B35").Value = "阅读"
$ws.Range("C35").Value = "为什么睡觉"
$ws.Range("D35").Value = 43

# 5. Append fresh date-only placeholder rows for the new day, ready for more entries
$ws.Range("A36").Value = 44949
$ws.Range("A37").Value = 44949
$ws.Range("A38").Value = 44949
$ws.Range("A39").Value = 44949
$ws.Range("A40").Value = 44949

# 6. Restore the view/selection state
$ws.Range("B36").Select()
